$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.6873773333333334
$ws.Range("H2").Value = 2.062132
$ws.Range("I2").Value = 0.02660947569874856
$ws.Range("J2").Value = 0.02660947569874856
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.900730666666667
$ws.Range("N2").Value = 17.702192
$ws.Range("O2").Value = 0.03970749001357476
$ws.Range("P2").Value = 0.03970749001357476
$ws.Range("Q2").Value = 4.056028510371556
$ws.Range("R2").Value = 36.504256593344
$ws.Range("S2").Value = 0.001056595490574519
$ws.Range("T2").Value = 0.001056595490574519
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.6873773333333334
$ws.Range("H3").Value = 2.062132
$ws.Range("I3").Value = 0.02660947569874856
$ws.Range("J3").Value = 0.02660947569874856
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 111.5917106666667
$ws.Range("N3").Value = 334.775132
$ws.Range("O3").Value = 0.7509284844884279
$ws.Range("P3").Value = 0.7509284844884279
$ws.Range("Q3").Value = 76.70561250015822
$ws.Range("R3").Value = 690.350512501424
$ws.Range("S3").Value = 0.01998181325949291
$ws.Range("T3").Value = 0.01998181325949291
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.6873773333333334
$ws.Range("H4").Value = 2.062132
$ws.Range("I4").Value = 0.02660947569874856
$ws.Range("J4").Value = 0.02660947569874856
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 31.11253633333333
$ws.Range("N4").Value = 93.337609
$ws.Range("O4").Value = 0.2093640254979974
$ws.Range("P4").Value = 0.2093640254979974
$ws.Range("Q4").Value = 21.38605225804312
$ws.Range("R4").Value = 192.474470322388
$ws.Range("S4").Value = 0.005571066948681137
$ws.Range("T4").Value = 0.005571066948681136
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.913984666666667
$ws.Range("H5").Value = 17.741954
$ws.Range("I5").Value = 0.2289398029860915
$ws.Range("J5").Value = 0.2289398029860915
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.900730666666667
$ws.Range("N5").Value = 17.702192
$ws.Range("O5").Value = 0.03970749001357476
$ws.Range("P5").Value = 0.03970749001357476
$ws.Range("Q5").Value = 34.89683068479645
$ws.Range("R5").Value = 314.071476163168
$ws.Range("S5").Value = 0.009090624940780001
$ws.Range("T5").Value = 0.009090624940779999
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 5.913984666666667
$ws.Range("H6").Value = 17.741954
$ws.Range("I6").Value = 0.2289398029860915
$ws.Range("J6").Value = 0.2289398029860915
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 111.5917106666667
$ws.Range("N6").Value = 334.775132
$ws.Range("O6").Value = 0.7509284844884279
$ws.Range("P6").Value = 0.7509284844884279
$ws.Range("Q6").Value = 659.9516658097698
$ws.Range("R6").Value = 5939.564992287927
$ws.Range("S6").Value = 0.1719174192954249
$ws.Range("T6").Value = 0.1719174192954249
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 5.913984666666667
$ws.Range("H7").Value = 17.741954
$ws.Range("I7").Value = 0.2289398029860915
$ws.Range("J7").Value = 0.2289398029860915
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 31.11253633333333
$ws.Range("N7").Value = 93.337609
$ws.Range("O7").Value = 0.2093640254979974
$ws.Range("P7").Value = 0.2093640254979974
$ws.Range("Q7").Value = 183.9990628164429
$ws.Range("R7").Value = 1655.991565347986
$ws.Range("S7").Value = 0.04793175874988656
$ws.Range("T7").Value = 0.04793175874988656
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 19.230689
$ws.Range("H8").Value = 57.692067
$ws.Range("I8").Value = 0.7444507213151601
$ws.Range("J8").Value = 0.7444507213151601
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.900730666666667
$ws.Range("N8").Value = 17.702192
$ws.Range("O8").Value = 0.03970749001357476
$ws.Range("P8").Value = 0.03970749001357476
$ws.Range("Q8").Value = 113.4751163234293
$ws.Range("R8").Value = 1021.276046910864
$ws.Range("S8").Value = 0.02956026958222025
$ws.Range("T8").Value = 0.02956026958222024
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 19.230689
$ws.Range("H9").Value = 57.692067
$ws.Range("I9").Value = 0.7444507213151601
$ws.Range("J9").Value = 0.7444507213151601
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 111.5917106666667
$ws.Range("N9").Value = 334.775132
$ws.Range("O9").Value = 0.7509284844884279
$ws.Range("P9").Value = 0.7509284844884279
$ws.Range("Q9").Value = 2145.98548280865
$ws.Range("R9").Value = 19313.86934527784
$ws.Range("S9").Value = 0.5590292519335102
$ws.Range("T9").Value = 0.5590292519335102
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 19.230689
$ws.Range("H10").Value = 57.692067
$ws.Range("I10").Value = 0.7444507213151601
$ws.Range("J10").Value = 0.7444507213151601
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.11253633333333
$ws.Range("N10").Value = 93.337609
$ws.Range("O10").Value = 0.2093640254979974
$ws.Range("P10").Value = 0.2093640254979974
$ws.Range("Q10").Value = 598.3155102275338
$ws.Range("R10").Value = 5384.839592047803
$ws.Range("S10").Value = 0.1558611997994298
$ws.Range("T10").Value = 0.1558611997994298